$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Sending cluster "MuSCs" -> Target cluster "ECs" (Ccl17/Ccr4 unchanged) ---
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Ccl17"
$ws.Range("C2").Value = "Ccr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.2468816666666667
$ws.Range("H2").Value = 0.740645
$ws.Range("I2").Value = 0.6299860588115711
$ws.Range("J2").Value = 0.6299860588115711
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04581866666666667
$ws.Range("N2").Value = 0.137456
$ws.Range("O2").Value = 0.4376255591461182
$ws.Range("P2").Value = 0.4376255591461182
$ws.Range("Q2").Value = 0.01131178879111111
$ws.Range("R2").Value = 0.10180609912
$ws.Range("S2").Value = 0.2756980012416732
$ws.Range("T2").Value = 0.2756980012416731

# --- Row 3: Sending cluster "MuSCs" -> Target cluster "FAPs" (Ccl17/Ccr4 unchanged) ---
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Ccl17"
$ws.Range("C3").Value = "Ccr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.2468816666666667
$ws.Range("H3").Value = 0.740645
$ws.Range("I3").Value = 0.6299860588115711
$ws.Range("J3").Value = 0.6299860588115711
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.05887966666666666
$ws.Range("N3").Value = 0.176639
$ws.Range("O3").Value = 0.5623744408538818
$ws.Range("P3").Value = 0.5623744408538818
$ws.Range("Q3").Value = 0.01453631023944444
$ws.Range("R3").Value = 0.130826792155
$ws.Range("S3").Value = 0.354288057569898
$ws.Range("T3").Value = 0.354288057569898

# --- Row 4 (new): Sending cluster "Resolving-Mac" -> Target cluster "ECs" ---
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("B4").Value = "Ccl17"
$ws.Range("C4").Value = "Ccr4"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1450026666666667
$ws.Range("H4").Value = 0.435008
$ws.Range("I4").Value = 0.3700139411884289
$ws.Range("J4").Value = 0.3700139411884289
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04581866666666667
$ws.Range("N4").Value = 0.137456
$ws.Range("O4").Value = 0.4376255591461182
$ws.Range("P4").Value = 0.4376255591461182
$ws.Range("Q4").Value = 0.006643828849777778
$ws.Range("R4").Value = 0.059794459648
$ws.Range("S4").Value = 0.1619275579044451
$ws.Range("T4").Value = 0.1619275579044451

# --- Row 5 (new): Sending cluster "Resolving-Mac" -> Target cluster "FAPs" ---
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Ccl17"
$ws.Range("C5").Value = "Ccr4"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1450026666666667
$ws.Range("H5").Value = 0.435008
$ws.Range("I5").Value = 0.3700139411884289
$ws.Range("J5").Value = 0.3700139411884289
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.05887966666666666
$ws.Range("N5").Value = 0.176639
$ws.Range("O5").Value = 0.5623744408538818
$ws.Range("P5").Value = 0.5623744408538818
$ws.Range("Q5").Value = 0.008537708679111111
$ws.Range("R5").Value = 0.07683937811199999
$ws.Range("S5").Value = 0.2080863832839838
$ws.Range("T5").Value = 0.2080863832839838
